$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update leg measurement data (D5, G5, D8) - the rest (H/K columns) are
# formula-driven and will recalculate automatically.
$ws.Range("D5").Value = 1580
$ws.Range("G5").Value = 2050
$ws.Range("D8").Value = 1250

# Update the active cell / selection on the sheet
$ws.Activate() | Out-Null
$ws.Range("F14").Select() | Out-Null

# Update the workbook window position (best effort; window chrome
# geometry is cosmetic and may not be persisted by every host).
$excel.ActiveWindow.Left = 5295
$excel.ActiveWindow.Top = 1170

$wb.Save() | Out-Null
